$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---

# Isa: extend the packs list
$ws.Range("C2").Value = "1 prima, 4 decouvertes, 4 packs 6, 4"

# Danny: fix punctuation + extend order, update amount paid
$ws.Range("C3").Value = "2 prima, 1 mix, 1 equilux, 1 pumpkin, 3 smarco, 2 ss, 1 chipmunk"
$ws.Range("D3").Value = 46

# Paco / Benoit: normalize "SanMarco" -> "smarco"
$ws.Range("C14").Value = "1 san marco, 1 antico, 5 smarco"
$ws.Range("C15").Value = "1 equilux, 5 smarco"

# --- Row 16 becomes a new customer "Theo B" (Clem's old data moves to row 17) ---
$ws.Range("A16").Value = "Theo B"
$ws.Range("C16").Value = "1 equilux, 1 antico"
$ws.Range("D16").Value = 15

# --- Row 17: Clem, re-added below Theo B with normalized packs text ---
$ws.Range("A17").Value = "Clem"
$ws.Range("C17").Value = "1 prima, 1 antico, 5 smarco"
$ws.Range("D17").Value = 51
$ws.Range("D17").NumberFormat = $ws.Range("D16").NumberFormat

# Match the selection recorded in the saved workbook
$ws.Range("G12").Select()
